$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 885.619
$ws.Range("J17").Value = 885.619
$ws.Range("L17").Value = 2656.857
$ws.Range("N17").Value = -2992.857
$ws.Range("H42").Value = 2778204
$ws.Range("J42").Value = 691.75
$ws.Range("L42").Value = 2075.25
$ws.Range("N42").Value = -2535.25
$ws.Range("H43").Value = 2052.9285
$ws.Range("J43").Value = 1650
$ws.Range("L43").Value = 1650
$ws.Range("N43").Value = -1788
$ws.Range("H46").Value = 1030.2693
$ws.Range("J46").Value = 1030.2693
$ws.Range("L46").Value = 3090.8079
$ws.Range("N46").Value = -3328.8079
$ws.Range("H60").Value = 1030.2693
$ws.Range("J60").Value = 1030.2693
$ws.Range("L60").Value = 3090.8079
$ws.Range("N60").Value = -4058.8079
$ws.Range("H64").Value = 39917.742
$ws.Range("J64").Value = 3003.4348
$ws.Range("L64").Value = 3003.4348
$ws.Range("N64").Value = -3499.4348
$ws.Range("H67").Value = 39917.742
$ws.Range("J67").Value = 3003.4348
$ws.Range("L67").Value = 3003.4348
$ws.Range("N67").Value = -4719.4348
$ws.Range("H86").Value = 2690.6775
$ws.Range("I86").Value = 1378.5625
$ws.Range("J86").Value = 4090.2666
$ws.Range("K86").Value = 1378.5625
$ws.Range("L86").Value = 4090.2666
$ws.Range("M86").Value = -255.5625
$ws.Range("N86").Value = -6336.2666
$ws.Range("H89").Value = 2690.6775
$ws.Range("I89").Value = 1378.5625
$ws.Range("J89").Value = 4090.2666
$ws.Range("K89").Value = 6892.8125
$ws.Range("L89").Value = 20451.333
$ws.Range("M89").Value = -1276.8125
$ws.Range("N89").Value = -31683.333
$ws.Range("H132").Value = 9267478
$ws.Range("I132").Value = 11372973
$ws.Range("K132").Value = 34118919
$ws.Range("M132").Value = -34116389

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24868.44
$ws.Range("I32").Value = 4934.3066
$ws.Range("J32").Value = 119938.92
$ws.Range("K32").Value = 4934.3066
$ws.Range("L32").Value = 119938.92
$ws.Range("M32").Value = -4647.3066
$ws.Range("N32").Value = -120512.92
$ws.Range("H122").Value = 3764.9167
$ws.Range("I122").Value = 3819.889
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 11459.667
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -9009.667000000001
$ws.Range("N122").Value = -15700
$ws.Range("H132").Value = 2067.4138
$ws.Range("I132").Value = 1544.8572
$ws.Range("J132").Value = 3439.125
$ws.Range("K132").Value = 4634.571599999999
$ws.Range("L132").Value = 10317.375
$ws.Range("M132").Value = -2104.571599999999
$ws.Range("N132").Value = -15377.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1031.4546
$ws.Range("J22").Value = 2251
$ws.Range("L22").Value = 2251
$ws.Range("N22").Value = -2951
$ws.Range("H31").Value = 25803.893
$ws.Range("I31").Value = 1007.19446
$ws.Range("J31").Value = 56586
$ws.Range("K31").Value = 1007.19446
$ws.Range("L31").Value = 56586
$ws.Range("M31").Value = -712.19446
$ws.Range("N31").Value = -57176
$ws.Range("H34").Value = 25803.893
$ws.Range("I34").Value = 1007.19446
$ws.Range("J34").Value = 56586
$ws.Range("K34").Value = 1007.19446
$ws.Range("L34").Value = 56586
$ws.Range("M34").Value = -805.19446
$ws.Range("N34").Value = -56990
$ws.Range("H62").Value = 2164.6667
$ws.Range("I62").Value = 1997.5
$ws.Range("K62").Value = 1997.5
$ws.Range("M62").Value = -1373.5
$ws.Range("H65").Value = 2164.6667
$ws.Range("I65").Value = 1997.5
$ws.Range("K65").Value = 9987.5
$ws.Range("M65").Value = -6867.5
$ws.Range("H107").Value = 1164.1818
$ws.Range("I107").Value = 1188.1111
$ws.Range("J107").Value = 1056.5
$ws.Range("K107").Value = 1188.1111
$ws.Range("L107").Value = 1056.5
$ws.Range("M107").Value = 731.8888999999999
$ws.Range("N107").Value = -4896.5
$ws.Range("H132").Value = 3183.5806
$ws.Range("I132").Value = 3042.0625
$ws.Range("J132").Value = 3334.5334
$ws.Range("K132").Value = 9126.1875
$ws.Range("L132").Value = 10003.6002
$ws.Range("M132").Value = -6596.1875
$ws.Range("N132").Value = -15063.6002
$ws.Range("H134").Value = 1403.5834
$ws.Range("I134").Value = 1311.4445
$ws.Range("J134").Value = 1680
$ws.Range("K134").Value = 3934.3335
$ws.Range("L134").Value = 5040
$ws.Range("M134").Value = -1399.3335
$ws.Range("N134").Value = -10110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2733.1667
$ws.Range("I62").Value = 799.6667
$ws.Range("J62").Value = 4666.6665
$ws.Range("K62").Value = 2399.0001
$ws.Range("L62").Value = 13999.9995
$ws.Range("M62").Value = -1713.0001
$ws.Range("N62").Value = -15371.9995
$ws.Range("H65").Value = 2733.1667
$ws.Range("I65").Value = 799.6667
$ws.Range("J65").Value = 4666.6665
$ws.Range("K65").Value = 7197.0003
$ws.Range("L65").Value = 41999.9985
$ws.Range("M65").Value = -3765.0003
$ws.Range("N65").Value = -48863.9985
$ws.Range("H131").Value = 799.36
$ws.Range("J131").Value = 824.21277
$ws.Range("L131").Value = 2472.63831
$ws.Range("N131").Value = -12552.63831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 64006.5
$ws.Range("I70").Value = 104190.35
$ws.Range("J70").Value = 6601
$ws.Range("K70").Value = 104190.35
$ws.Range("L70").Value = 6601
$ws.Range("M70").Value = -103920.35
$ws.Range("N70").Value = -7141
$ws.Range("H73").Value = 64006.5
$ws.Range("I73").Value = 104190.35
$ws.Range("J73").Value = 6601
$ws.Range("K73").Value = 104190.35
$ws.Range("L73").Value = 6601
$ws.Range("M73").Value = -103254.35
$ws.Range("N73").Value = -8473
$ws.Range("H113").Value = 1564.6389
$ws.Range("I113").Value = 1390.1875
$ws.Range("J113").Value = 1704.2
$ws.Range("K113").Value = 1390.1875
$ws.Range("L113").Value = 1704.2
$ws.Range("M113").Value = 779.8125
$ws.Range("N113").Value = -6044.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1626.2
$ws.Range("I136").Value = 1542.2
$ws.Range("J136").Value = 1878.2
$ws.Range("K136").Value = 4626.6
$ws.Range("L136").Value = 5634.6
$ws.Range("M136").Value = -2076.6
$ws.Range("N136").Value = -10734.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 25792.715
$ws.Range("J42").Value = 25792.715
$ws.Range("L42").Value = 25792.715
$ws.Range("N42").Value = -26548.715
$ws.Range("H43").Value = 52000
$ws.Range("I43").Value = 52000
$ws.Range("K43").Value = 52000
$ws.Range("M43").Value = -51851
